$wb = $excel.ActiveWorkbook

# ---- Sheet "DataSet" (sheet2.xml): Pass/Fail result column P ----
$ws2 = $wb.Worksheets.Item("DataSet")

# Cells P4, P7, P10, P15, P18 already carry a "quote-prefixed text" style
# (quotePrefix="1"), so re-assigning with a leading apostrophe keeps that
# style index unchanged while turning the cell into a literal text value.
$ws2.Range("P4").Value = "'Passed"
$ws2.Range("P7").Value = "'Passed"
$ws2.Range("P10").Value = "'passed"
$ws2.Range("P15").Value = "'passed"
$ws2.Range("P18").Value = "'passed"

# P22 uses a plain (non quote-prefixed) style, so assign the literal text
# directly -- no apostrophe, so Excel doesn't switch it to a quote-prefix
# style.
$ws2.Range("P22").Value = "passed"

# ---- Sheet "NegativeTests" (sheet3.xml): Pass/Fail result column J ----
$ws3 = $wb.Worksheets.Item("NegativeTests")

# J2..J5 already have a quote-prefixed style -- same trick as above.
$ws3.Range("J2").Value = "'Passed"
$ws3.Range("J3").Value = "'Passed"
$ws3.Range("J4").Value = "'Passed"
$ws3.Range("J5").Value = "'Passed"

# J6..J8 are brand-new cells. Pull the quote-prefixed / centered style used
# by the rest of the "Passed" column (e.g. DataSet!P4) onto them first, then
# write the value.
$ws2.Range("P4").Copy()
$ws3.Range("J6").PasteSpecial(-4122)
$ws3.Range("J7").PasteSpecial(-4122)
$ws3.Range("J8").PasteSpecial(-4122)

$ws3.Range("J6").Value = "'Passed"
$ws3.Range("J7").Value = "'Passed"
$ws3.Range("J8").Value = "'Passed"
